# Fruta / hortaliza, semanal
# The weekly refresh reshuffles which dated record occupies which row of
# the "Breva" subset for Vega Central Mapocho de Santiago. Only the
# per-record fields (Fecha, Calidad, Volumen, Precio minimo/maximo/
# promedio ponderado, Unidad de comercializacion, Origen, Precio $/Kg and
# Kg / unidad) move between rows 2-5 and 8-13; the descriptive columns
# (Mercado, Region, Producto, Categoria, Variedad, etc.) stay put since
# they are identical for every row already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44553
$ws.Range("L2").Value = 'Especial'
$ws.Range("N2").Value = 22000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 22000
$ws.Range("S2").Value = 3667

# Row 3
$ws.Range("D3").Value = 44553
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("R3").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S3").Value = 3000

# Row 4
$ws.Range("D4").Value = 44561
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = '$/bandeja 6 kilos'
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 6

# Row 5
$ws.Range("D5").Value = 44204
$ws.Range("M5").Value = 110
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7500
$ws.Range("P5").Value = 7318
$ws.Range("S5").Value = 1045

# Row 8
$ws.Range("D8").Value = 44550
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 24000
$ws.Range("Q8").Value = '$/bandeja 7 kilos'
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 3429
$ws.Range("T8").Value = 7

# Row 9
$ws.Range("D9").Value = 44558
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 22000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 22000
$ws.Range("S9").Value = 3667

# Row 10
$ws.Range("D10").Value = 44558
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("Q10").Value = '$/bandeja 6 kilos'
$ws.Range("S10").Value = 3000
$ws.Range("T10").Value = 6

# Row 11
$ws.Range("D11").Value = 44187
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 45
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("R11").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S11").Value = 2000

# Row 12
$ws.Range("D12").Value = 44187
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("Q12").Value = '$/bandeja 7 kilos'
$ws.Range("S12").Value = 1714
$ws.Range("T12").Value = 7

# Row 13
$ws.Range("D13").Value = 44572
$ws.Range("M13").Value = 65
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 3333
